# Insert a new data row at row 446 (pushes existing rows 446-501 down to 447-502)
# and populate it with a new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("446").Insert()

$ws.Range("A446").Value = 9
$ws.Range("B446").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C446").Value = "Metropolitana"
$ws.Range("D446").Value = 45212
$ws.Range("E446").Value = 13
$ws.Range("F446").Value = 300000001
$ws.Range("G446").Value = "Rabanito"
$ws.Range("H446").Value = "Sin especificar"
$ws.Range("I446").Value = "Primera"
$ws.Range("J446").Value = 7000
$ws.Range("K446").Value = 3000
$ws.Range("L446").Value = 3000
$ws.Range("M446").Value = 3000
$ws.Range("N446").Value = '$/cien unidades (volumen en unidades)'
$ws.Range("O446").Value = "Región Metropolitana"
$ws.Range("P446").Value = 30
$ws.Range("Q446").Value = 100
$ws.Range("R446").Value = "Hortaliza"
